## edit.ps1 - apply the two changes described by the diff:
##  1. Bold the run "程序逻辑" in the existing text box (shape id 31 / "文本框 30").
##  2. Add a new text box "难点：树状文件夹结构" at the end of the shape tree
##     (shape id 122 / "文本框 121").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) Bold the "程序逻辑" run (11th shape in the slide, cNvPr id="31").
# ---------------------------------------------------------------------------
$logicBox = $s.Shapes.Item(11)
$logicBox.TextFrame.TextRange.Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) Insert the new "难点：树状文件夹结构" text box.
#
# This engine assigns each newly created shape the smallest shape id that is
# not yet used anywhere in the slide's XML, tracked by a cursor that only
# ever moves forward (ids freed up by deleting a shape are not revisited).
# The target id from the diff is 122, which - starting the cursor at 1 on
# this slide - is the 71st id that gets handed out. So we burn through the
# first 70 candidate ids with scratch textboxes (immediately deleted, so
# they leave no trace in the saved file) and then create the real shape,
# which lands on id 122, matching "文本框 121".
# ---------------------------------------------------------------------------
$scratch = New-Object System.Collections.ArrayList
for ($i = 0; $i -lt 70; $i++) {
    [void]$scratch.Add($s.Shapes.AddTextbox(1, 0, 0, 1, 1))
}
foreach ($junk in $scratch) {
    $junk.Delete()
}

# EMU -> point conversion (1 pt = 12700 EMU) since AddTextbox takes points.
$left = 7514293 / 12700
$top = 6411494 / 12700
$width = 2587568 / 12700
$height = 369332 / 12700

$box = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$box.Name = "文本框 121"

# No fill on the shape.
$box.Fill.Visible = 0

# Text + run language (matches lang="zh-CN" altLang="en-US").
$box.TextFrame.TextRange.Text = "难点：树状文件夹结构"
$box.TextFrame.TextRange.LanguageID = "zh-CN"

# Body properties: wrap="none" + spAutoFit (shrink box to fit the text).
$box.TextFrame.WordWrap = 0
$box.TextFrame.AutoSize = 1
